# Commit: "Rename arc to link"
#
# Rename the "arcs" worksheet to "links" (the tab stays in the same
# position / keeps its sheetId, only the display name changes), and
# restore the sheet's last-known selection.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("arcs")

# Select the sheet and set the cell selection that was recorded for it
# before renaming (activating it mirrors the original file, where this
# sheet is the tabSelected / active one).
$ws.Activate()
$ws.Range("Q17").Select()

$ws.Name = "links"
